# Change the default font of the "content placeholder" (内容占位符, idx=1)
# text box on each of the three "Title and Content" slide layouts to
# Times New Roman, covering all five outline levels that live inside it.
#
# These are the layouts wired up to the slide master as:
#   CustomLayouts.Item(2) -> slideLayout2.xml ("标题和内容")
#   CustomLayouts.Item(3) -> slideLayout3.xml ("1_标题和内容")
#   CustomLayouts.Item(4) -> slideLayout4.xml ("2_标题和内容")
# and in every one of them shape #1 is the multi-level content
# placeholder whose paragraphs (lvl 0-4) should render in Times New Roman.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster

foreach ($layoutIndex in 2, 3, 4) {
    $layout = $master.CustomLayouts.Item($layoutIndex)
    $shape = $layout.Shapes.Item(1)

    if ($shape.HasTextFrame) {
        $textRange = $shape.TextFrame.TextRange
        $textRange.Font.Name = "Times New Roman"
    }
}
